$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (rows 26-27), pushing the
# existing rows 26-48 down to 28-50 (this is a weekly data refresh: the two
# newest records are added at the top of the Membrillo / Talca subset).
$ws.Rows("26:27").Insert()

# New row 26 - Especial
$ws.Range("A26").Value = 5
$ws.Range("B26").Value = "Macroferia Regional de Talca"
$ws.Range("C26").Value = "Maule"
$ws.Range("D26").Value = 44741
$ws.Range("E26").Value = 7
$ws.Range("F26").Value = "Fruta"
$ws.Range("G26").Value = 100104
$ws.Range("H26").Value = "Frutos de pepita"
$ws.Range("I26").Value = 100104003
$ws.Range("J26").Value = "Membrillo"
$ws.Range("K26").Value = "Champion"
$ws.Range("L26").Value = "Especial"
$ws.Range("M26").Value = 180
$ws.Range("N26").Value = 12000
$ws.Range("O26").Value = 12000
$ws.Range("P26").Value = 12000
$ws.Range("Q26").Value = "$/caja 18 kilos granel"
$ws.Range("R26").Value = "Región de O'Higgins"
$ws.Range("S26").Value = 667
$ws.Range("T26").Value = 18

# New row 27 - Primera
$ws.Range("A27").Value = 5
$ws.Range("B27").Value = "Macroferia Regional de Talca"
$ws.Range("C27").Value = "Maule"
$ws.Range("D27").Value = 44741
$ws.Range("E27").Value = 7
$ws.Range("F27").Value = "Fruta"
$ws.Range("G27").Value = 100104
$ws.Range("H27").Value = "Frutos de pepita"
$ws.Range("I27").Value = 100104003
$ws.Range("J27").Value = "Membrillo"
$ws.Range("K27").Value = "Champion"
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 230
$ws.Range("N27").Value = 10000
$ws.Range("O27").Value = 10000
$ws.Range("P27").Value = 10000
$ws.Range("Q27").Value = "$/caja 18 kilos granel"
$ws.Range("R27").Value = "Región de O'Higgins"
$ws.Range("S27").Value = 556
$ws.Range("T27").Value = 18
